$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume 1h change (E) columns for every coin row (2-51).
# Force the D-column cell to Text format before writing so values such as
# "1.000" or "0.9984" are stored as literal strings (matching the source data,
# which is type inlineStr) instead of being auto-converted to numbers by
# Excel value-assignment type inference. The style is reset to "Normal" right
# after the write so the temporary Text number-format does not linger on the
# cell (the original cells carry no explicit style).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.854.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.632.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4710'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '38.80'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2552'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06061'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06938'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.636.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6023'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.311'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '72.78'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9998'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '24.856.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006544'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.846.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.337'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.533'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.200'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '132.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.373'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '103.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.622'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.766'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07710'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.516'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9991'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.04270'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.583'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9170'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5761'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.533'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01533'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9991'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8089'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.03'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.759'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3671'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.694'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05218'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1085'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.987'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '29.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.26%  '
